# "made chagnes for month append" — roll the inventory/sales report from
# April-2024 to May-2024: bump the report-month date in B1 and refresh every
# Outwards Qty / Rate / Amount figure (and the summary block below it) with
# the new month's numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 — report month (date serial): Apr-24 -> May-24
$ws.Range("B1").Value = 45413

# Row 4 — MCF
$ws.Range("B4").Value = 53859.21
$ws.Range("C4").Value = 260.95243506171
$ws.Range("D4").Value = 14054692

# Row 5 — WMF
$ws.Range("B5").Value = 41654.76
$ws.Range("C5").Value = 261.6926180825432
$ws.Range("D5").Value = 10900743.2

# Row 7 — MONOFILAMENT FABRIC INSECT NET
$ws.Range("B7").Value = 135.58
$ws.Range("C7").Value = 497.8610414515415
$ws.Range("D7").Value = 67500

# Row 8 — MONOFILAMENT FABRIC HAPPA (now has Qty/Amount too, previously blank)
$ws.Range("B8").Value = 10.02
$ws.Range("C8").Value = 673.6526946107784
$ws.Range("D8").Value = 6750

# Row 9 — NWF/Yarn
$ws.Range("B9").Value = 9050
$ws.Range("C9").Value = 131.50552486187846
$ws.Range("D9").Value = 1190125

# Row 10 — MSN
$ws.Range("B10").Value = 104709.57
$ws.Range("C10").Value = 250.40509859796003
$ws.Range("D10").Value = 26219810.2

# Row 11 — TSN
$ws.Range("B11").Value = 42.9
$ws.Range("C11").Value = 242.42424242424244
$ws.Range("D11").Value = 10400

# Row 12 — PP Woven Sacks
$ws.Range("B12").Value = 298.25
$ws.Range("C12").Value = 237.0829840737636
$ws.Range("D12").Value = 70710

# Row 13 — ANTI BIRD NET / Rope/MULCH/FIBC
$ws.Range("B13").Value = 6605
$ws.Range("C13").Value = 148.13641180923543
$ws.Range("D13").Value = 978441

# Row 15 — Weed Mat 1.25 Mtrs Black
$ws.Range("B15").Value = 1192.03
$ws.Range("C15").Value = 256.7492428881823
$ws.Range("D15").Value = 306052.8

# Row 16 — Flexible Intermediate Bulk Container
$ws.Range("B16").Value = 949.42
$ws.Range("C16").Value = 193.1788881632997
$ws.Range("D16").Value = 183407.9

# Row 17 — Packing Materials / Old use Batteries
$ws.Range("B17").Value = 9087.599999999999
$ws.Range("C17").Value = 170.4533320128527
$ws.Range("D17").Value = 1549011.7

# Row 20 — HDPE Monofilament Waste
$ws.Range("B20").Value = 4134
$ws.Range("D20").Value = 41340

# Row 22 — Raw Material
$ws.Range("B22").Value = 54225
$ws.Range("C22").Value = 96.55002305209774
$ws.Range("D22").Value = 5235425

# Row 23 — (unlabeled continuation row)
$ws.Range("B23").Value = 58564
$ws.Range("C23").Value = 90.87647701659722
$ws.Range("D23").Value = 5322090

# Row 24 — Grand total
$ws.Range("B24").Value = 172361.17
$ws.Range("C24").Value = 191.9858858001486
$ws.Range("D24").Value = 33090911.9

# Row 25 — Other Income (now carries a Qty figure too)
$ws.Range("B25").Value = 172361.17
$ws.Range("D25").Value = 10.94

# Row 26 — (unlabeled subtotal row, now carries a Qty figure too)
$ws.Range("B26").Value = 0
$ws.Range("D26").Value = 33049572

# Row 27 — Gross sales
$ws.Range("D27").Value = 35651792.9

# Row 28 — tax
$ws.Range("D28").Value = 2566561

# Row 29 — TCS (amount cleared entirely — no longer applicable this month)
$ws.Range("D29").ClearContents()

# Row 30 — (unlabeled subtotal row)
$ws.Range("D30").Value = 2566561

# Row 31 — (unlabeled subtotal row)
$ws.Range("D31").Value = 33085231.9

# Row 32 — Discount
$ws.Range("D32").Value = 101080

# Row 33 — (unlabeled subtotal row)
$ws.Range("D33").Value = 33090911.9

# Row 34 — Credit Note
$ws.Range("D34").Value = 95400

# Row 36 — PAL I- FINAL SALES +Less Waste & Discount
$ws.Range("D36").Value = 32948492

# Row 39 — RM Purchase for sales
$ws.Range("B39").Value = 54225
$ws.Range("C39").Value = 90.36
$ws.Range("D39").Value = 4899771
